$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enhancement #24 : Collected company details
# Append a new company-details row (row 45) duplicating the
# IBCS-PRIMAX Software (Bangladesh) Limited entry (same data as row 10),
# complete with its hyperlinks, matching row height and cell formatting.

$ws.Range("A45").Value = "IBCS-PRIMAX Software(Bangladesh) Limited"
$ws.Range("A45").VerticalAlignment = -4160

$ws.Range("B45").Value = "House # 51, Road # 10A, Dhanmondi R/A,`nDhaka-1209, Bangladesh"
$ws.Range("B45").WrapText = $true

$ws.Range("C45").Value = "info@ibcs-primax.com"
$ws.Hyperlinks.Add($ws.Range("C45"), "mailto:info@ibcs-primax.com")
$ws.Range("C45").VerticalAlignment = -4160

$ws.Range("D45").Value = 58153970
$ws.Range("D45").HorizontalAlignment = -4131
$ws.Range("D45").VerticalAlignment = -4160

$ws.Range("E45").Value = "http://www.ibcs-primax.com"
$ws.Hyperlinks.Add($ws.Range("E45"), "http://www.ibcs-primax.com/")
$ws.Range("E45").VerticalAlignment = -4160

$ws.Rows("45").RowHeight = 45

# Move the selection down past the newly-added row, and scroll the
# window so the new row is visible near the bottom of the view.
[void]$ws.Range("A46").Select()
$excel.ActiveWindow.ScrollRow = 41
